# Add DOI/citations data to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "doi"
$ws.Range("B1").Value = "citations"

# Data rows
$ws.Range("A2").Value = 10.276783
$ws.Range("B2").Value = 66

$ws.Range("A3").Value = 10.264275697
$ws.Range("B3").Value = 1

# Widen the "doi" column to fit its contents
$ws.Columns.Item(1).ColumnWidth = 17.28515625

# Leave the selection where it would land after typing the last value
# and pressing Enter (one row below the data, in column B)
$ws.Range("B4").Select() | Out-Null
